$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a text value into a cell while preventing Excel's automatic
# conversion of date/number-looking strings (e.g. "2026-02-17") into a date
# serial value. The cell's NumberFormat is forced to Text before the write,
# then reset back to the default "Normal" style so no stray formatting is
# left behind.
# ---------------------------------------------------------------------------
function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------------
# Helper: force a cell to exist but stay blank (matches source rows that
# contain an explicit-but-empty cell, e.g. an open trade with no exit price
# / exit reason yet).
# ---------------------------------------------------------------------------
function Set-BlankCell($cell) {
    $cell.NumberFormat = "@"
    $cell.Value = ""
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------------
# Sheet "Summary" updates
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.49
$summary.Range("B4").Value = -0.51
$summary.Range("B6").Value = 186
$summary.Range("B8").Value = 78
$summary.Range("B9").Value = 40.86

# ---------------------------------------------------------------------------
# Sheet "Strategy Status" updates (volatility_scorer row, row 12)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C12").Value = 99.44
$status.Range("D12").Value = 11
$status.Range("E12").Value = -0.56
$status.Range("F12").Value = -0.56
$status.Range("G12").Value = 36.36

# ---------------------------------------------------------------------------
# Sheet "All Trades" - append trade #186 (row 187) and trade #187 (row 188)
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

Set-TextValue $allTrades.Cells.Item(187, 2) "2026-02-17"
Set-TextValue $allTrades.Cells.Item(187, 3) "10:07:05"
$allTrades.Cells.Item(187, 1).Value = 186
$allTrades.Cells.Item(187, 4).Value = "volatility_scorer"
$allTrades.Cells.Item(187, 5).Value = "NEUTRAL"
$allTrades.Cells.Item(187, 6).Value = 0.37
$allTrades.Cells.Item(187, 7).Value = 0.29
$allTrades.Cells.Item(187, 8).Value = "CLOSED"
$allTrades.Cells.Item(187, 9).Value = -21.6216
$allTrades.Cells.Item(187, 10).Value = -0.08
$allTrades.Cells.Item(187, 11).Value = 99.44
$allTrades.Cells.Item(187, 12).Value = 0
$allTrades.Cells.Item(187, 13).Value = 0
$allTrades.Cells.Item(187, 14).Value = 0.85
$allTrades.Cells.Item(187, 15).Value = "Low vol market (score: inf) - ideal for market making"
$allTrades.Cells.Item(187, 16).Value = "early_exit"
$allTrades.Cells.Item(187, 17).Value = 0.17

Set-TextValue $allTrades.Cells.Item(188, 2) "2026-02-17"
Set-TextValue $allTrades.Cells.Item(188, 3) "10:07:06"
$allTrades.Cells.Item(188, 1).Value = 187
$allTrades.Cells.Item(188, 4).Value = "MarketMaking"
$allTrades.Cells.Item(188, 5).Value = "UP"
$allTrades.Cells.Item(188, 6).Value = 0.63
Set-BlankCell $allTrades.Cells.Item(188, 7)
$allTrades.Cells.Item(188, 8).Value = "OPEN"
$allTrades.Cells.Item(188, 9).Value = 0
$allTrades.Cells.Item(188, 10).Value = 0
$allTrades.Cells.Item(188, 11).Value = 100.0480687506789
$allTrades.Cells.Item(188, 12).Value = 0
$allTrades.Cells.Item(188, 13).Value = 0
$allTrades.Cells.Item(188, 14).Value = 0.6
$allTrades.Cells.Item(188, 15).Value = "Normal spread capture: 19600 bps"
Set-BlankCell $allTrades.Cells.Item(188, 16)
$allTrades.Cells.Item(188, 17).Value = 0

# ---------------------------------------------------------------------------
# Sheet "volatility_scorer" - append trade #186 (row 12)
# ---------------------------------------------------------------------------
$volScorer = $wb.Worksheets.Item("volatility_scorer")

Set-TextValue $volScorer.Cells.Item(12, 2) "2026-02-17"
Set-TextValue $volScorer.Cells.Item(12, 3) "10:07:05"
$volScorer.Cells.Item(12, 1).Value = 186
$volScorer.Cells.Item(12, 4).Value = "volatility_scorer"
$volScorer.Cells.Item(12, 5).Value = "NEUTRAL"
$volScorer.Cells.Item(12, 6).Value = 0.37
$volScorer.Cells.Item(12, 7).Value = 0.29
$volScorer.Cells.Item(12, 8).Value = "CLOSED"
$volScorer.Cells.Item(12, 9).Value = -21.6216
$volScorer.Cells.Item(12, 10).Value = -0.08
$volScorer.Cells.Item(12, 11).Value = 99.44
$volScorer.Cells.Item(12, 12).Value = 0
$volScorer.Cells.Item(12, 13).Value = 0
$volScorer.Cells.Item(12, 14).Value = 0.85
$volScorer.Cells.Item(12, 15).Value = "Low vol market (score: inf) - ideal for market making"
$volScorer.Cells.Item(12, 16).Value = "early_exit"
$volScorer.Cells.Item(12, 17).Value = 0.17

# ---------------------------------------------------------------------------
# Sheet "MarketMaking" - append trade #187 (row 177)
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")

Set-TextValue $marketMaking.Cells.Item(177, 2) "2026-02-17"
Set-TextValue $marketMaking.Cells.Item(177, 3) "10:07:06"
$marketMaking.Cells.Item(177, 1).Value = 187
$marketMaking.Cells.Item(177, 4).Value = "MarketMaking"
$marketMaking.Cells.Item(177, 5).Value = "UP"
$marketMaking.Cells.Item(177, 6).Value = 0.63
Set-BlankCell $marketMaking.Cells.Item(177, 7)
$marketMaking.Cells.Item(177, 8).Value = "OPEN"
$marketMaking.Cells.Item(177, 9).Value = 0
$marketMaking.Cells.Item(177, 10).Value = 0
$marketMaking.Cells.Item(177, 11).Value = 100.0480687506789
$marketMaking.Cells.Item(177, 12).Value = 0
$marketMaking.Cells.Item(177, 13).Value = 0
$marketMaking.Cells.Item(177, 14).Value = 0.6
$marketMaking.Cells.Item(177, 15).Value = "Normal spread capture: 19600 bps"
Set-BlankCell $marketMaking.Cells.Item(177, 16)
$marketMaking.Cells.Item(177, 17).Value = 0
